$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: I1 = "I0", J1 = "IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font, border, centered alignment) from the
# existing header cell H1 onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows 2-29: columns I (I0) and J (IF) ---
$data = @(
    @(2, 2, 3),
    @(3, 1, 2),
    @(4, 5, 6),
    @(5, 6, 7),
    @(6, 6, 8),
    @(7, 5, 7),
    @(8, 8, 8),
    @(9, 7, 7),
    @(10, 1, 1),
    @(11, 10, 10),
    @(12, 3, 3),
    @(13, 7, 7),
    @(14, 6, 8),
    @(15, 6, 7),
    @(16, 6, 7),
    @(17, 7, 7),
    @(18, 2, 3),
    @(19, 7, 8),
    @(20, 1, 2),
    @(21, 9, 9),
    @(22, 2, 4),
    @(23, 8, 8),
    @(24, 8, 9),
    @(25, 6, 7),
    @(26, 8, 9),
    @(27, 4, 4),
    @(28, 7, 7),
    @(29, 8, 9)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}

Write-Host "I0 and IF columns added"
